# Fill in the TSP distance-matrix diagonal (distance from a city to itself = 0).
# Row r (2..77) is missing its diagonal cell in the same-lettered column as the
# row number (row 2 -> B2, row 3 -> C3, ... row 77 -> BY77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$diagCells = @(
    "B2","C3","D4","E5","F6","G7","H8","I9","J10","K11","L12","M13","N14","O15",
    "P16","Q17","R18","S19","T20","U21","V22","W23","X24","Y25","Z26","AA27","AB28",
    "AC29","AD30","AE31","AF32","AG33","AH34","AI35","AJ36","AK37","AL38","AM39","AN40",
    "AO41","AP42","AQ43","AR44","AS45","AT46","AU47","AV48","AW49","AX50","AY51","AZ52",
    "BA53","BB54","BC55","BD56","BE57","BF58","BG59","BH60","BI61","BJ62","BK63","BL64",
    "BM65","BN66","BO67","BP68","BQ69","BR70","BS71","BT72","BU73","BV74","BW75","BX76",
    "BY77"
)

foreach ($cell in $diagCells) {
    $ws.Range($cell).Value = 0
}

# Leave the final selection on B1, matching the saved workbook state.
$ws.Range("B1").Select()
